# Insert a new data row at row 48 (pushing existing rows 48-136 down to
# rows 49-137) and populate it with a new "Albahaca" price record for
# "Terminal La Palmera de La Serena" (Coquimbo).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 48; this shifts row 48
# (and everything below it) down by one, to row 49, etc.
$ws.Rows.Item(48).Insert()

# The record that used to live in row 48 is now in row 49. Most columns
# of this dataset repeat the same values down the whole block (Mercado ID,
# Mercado, Region, Codreg, Categoria ID, Categoria, Variedad, Calidad,
# Unidad de comercializacion, Origen, Kg o Unidades, Clasificacion), so
# carry those over into the freshly inserted row 48 from row 49.
$ws.Range("A48").Value = $ws.Range("A49").Value2
$ws.Range("B48").Value = $ws.Range("B49").Value2
$ws.Range("C48").Value = $ws.Range("C49").Value2
$ws.Range("E48").Value = $ws.Range("E49").Value2
$ws.Range("F48").Value = $ws.Range("F49").Value2
$ws.Range("G48").Value = $ws.Range("G49").Value2
$ws.Range("H48").Value = $ws.Range("H49").Value2
$ws.Range("I48").Value = $ws.Range("I49").Value2
$ws.Range("N48").Value = $ws.Range("N49").Value2
$ws.Range("O48").Value = $ws.Range("O49").Value2
$ws.Range("Q48").Value = $ws.Range("Q49").Value2
$ws.Range("R48").Value = $ws.Range("R49").Value2

# Now write the new record's own values (the ones that actually differ).
$ws.Range("D48").Value = 44883
$ws.Range("J48").Value = 1000
$ws.Range("K48").Value = 3500
$ws.Range("L48").Value = 4000
$ws.Range("M48").Value = 3750
$ws.Range("P48").Value = 3750
